$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Update the "Situacao da OS" status text from "Aberta" to
#    "Plano Entregue".
# ------------------------------------------------------------------
$replaceRange = $d.Content
$replaceFind = $replaceRange.Find
$replaced = $replaceFind.Execute("Aberta", $true, $true, $false, $false, $false, $true, 1, $false, "Plano Entregue", 2)

# ------------------------------------------------------------------
# 2. Word keeps the "_GoBack" bookmark pinned to the location of the
#    most recent edit. Remove the old bookmark and recreate it right
#    after the text that was just replaced.
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# Locate the text we just inserted (it is the first occurrence of
# "Plano Entregue" in the document, since the original status field
# precedes the history table further down that already contained the
# phrase).
$locateRange = $d.Content
$locateFind = $locateRange.Find
$located = $locateFind.Execute("Plano Entregue", $true, $true, $false, $false, $false, $true, 0, $false)

if ($located) {
    $locateRange.Collapse(0)
    $d.Bookmarks.Add("_GoBack", $locateRange)
}
